# Updates the "cryptos" price list worksheet with refreshed price/volume
# figures (and, for a few re-ranked coins, updated name/link/price cells)
# as published by the upstream GitHub Actions data refresh job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry describes the single-cell update to apply. Columns D (Price)
# and E (Volume(1h)) hold values that look numeric/percentage-like
# (e.g. "1.000", "30.261.63", "  -0.45%  ") but must remain plain text,
# exactly as scraped from coinranking.com, so force a text number format
# on those columns before assigning the value.
$updates = @(
    @{ Cell = 'D2'; Value = '30.270.93' }
    @{ Cell = 'E2'; Value = '  -0.42%  ' }
    @{ Cell = 'D3'; Value = '1.857.14' }
    @{ Cell = 'E3'; Value = '  -1.12%  ' }
    @{ Cell = 'E4'; Value = '  +0.04%  ' }
    @{ Cell = 'D5'; Value = '232.87' }
    @{ Cell = 'E5'; Value = '  -2.40%  ' }
    @{ Cell = 'E6'; Value = '  +0.05%  ' }
    @{ Cell = 'D7'; Value = '0.4751' }
    @{ Cell = 'E7'; Value = '  -0.87%  ' }
    @{ Cell = 'D8'; Value = '0.2751' }
    @{ Cell = 'E8'; Value = '  -2.59%  ' }
    @{ Cell = 'D9'; Value = '0.06423' }
    @{ Cell = 'E9'; Value = '  -1.46%  ' }
    @{ Cell = 'D10'; Value = '1.862.70' }
    @{ Cell = 'E10'; Value = '  -0.83%  ' }
    @{ Cell = 'E11'; Value = '  -0.57%  ' }
    @{ Cell = 'D12'; Value = '16.09' }
    @{ Cell = 'E12'; Value = '  -3.47%  ' }
    @{ Cell = 'D13'; Value = '4.987' }
    @{ Cell = 'E13'; Value = '  -2.22%  ' }
    @{ Cell = 'D14'; Value = '85.31' }
    @{ Cell = 'E14'; Value = '  -3.32%  ' }
    @{ Cell = 'D15'; Value = '0.6326' }
    @{ Cell = 'E15'; Value = '  -4.32%  ' }
    @{ Cell = 'D16'; Value = '30.251.05' }
    @{ Cell = 'E16'; Value = '  -0.40%  ' }
    @{ Cell = 'D17'; Value = '1.000' }
    @{ Cell = 'E17'; Value = '  +0.05%  ' }
    @{ Cell = 'D18'; Value = '12.77' }
    @{ Cell = 'E18'; Value = '  -4.07%  ' }
    @{ Cell = 'D19'; Value = '0.000007317' }
    @{ Cell = 'E19'; Value = '  -3.85%  ' }
    @{ Cell = 'B20'; Value = 'BitcoinCash' }
    @{ Cell = 'C20'; Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch' }
    @{ Cell = 'D20'; Value = '224.74' }
    @{ Cell = 'E20'; Value = '  +2.40%  ' }
    @{ Cell = 'B21'; Value = 'WrappedliquidstakedEther2.0' }
    @{ Cell = 'C21'; Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth' }
    @{ Cell = 'D21'; Value = '2.091.64' }
    @{ Cell = 'E21'; Value = '  -1.13%  ' }
    @{ Cell = 'D22'; Value = '1.001' }
    @{ Cell = 'E22'; Value = '  +0.06%  ' }
    @{ Cell = 'D23'; Value = '5.095' }
    @{ Cell = 'E23'; Value = '  -3.93%  ' }
    @{ Cell = 'D24'; Value = '6.016' }
    @{ Cell = 'E24'; Value = '  -3.08%  ' }
    @{ Cell = 'D25'; Value = '167.40' }
    @{ Cell = 'E25'; Value = '  -0.09%  ' }
    @{ Cell = 'D26'; Value = '9.223' }
    @{ Cell = 'E26'; Value = '  -1.71%  ' }
    @{ Cell = 'D27'; Value = '17.80' }
    @{ Cell = 'E27'; Value = '  -3.63%  ' }
    @{ Cell = 'D28'; Value = '1.859' }
    @{ Cell = 'E28'; Value = '  -5.99%  ' }
    @{ Cell = 'D29'; Value = '0.1028' }
    @{ Cell = 'E29'; Value = '  +9.78%  ' }
    @{ Cell = 'D30'; Value = '1.380' }
    @{ Cell = 'E30'; Value = '  -5.58%  ' }
    @{ Cell = 'D31'; Value = '4.217' }
    @{ Cell = 'E31'; Value = '  -2.65%  ' }
    @{ Cell = 'D32'; Value = '3.898' }
    @{ Cell = 'E32'; Value = '  -3.39%  ' }
    @{ Cell = 'D33'; Value = '0.04884' }
    @{ Cell = 'E33'; Value = '  -3.16%  ' }
    @{ Cell = 'D34'; Value = '1.147' }
    @{ Cell = 'E34'; Value = '  -4.78%  ' }
    @{ Cell = 'D35'; Value = '0.7247' }
    @{ Cell = 'E35'; Value = '  -2.71%  ' }
    @{ Cell = 'D36'; Value = '0.9993' }
    @{ Cell = 'E36'; Value = '  +0.12%  ' }
    @{ Cell = 'D37'; Value = '2.684' }
    @{ Cell = 'E37'; Value = '  -0.87%  ' }
    @{ Cell = 'D38'; Value = '0.01902' }
    @{ Cell = 'E38'; Value = '  +4.06%  ' }
    @{ Cell = 'D39'; Value = '2.628' }
    @{ Cell = 'E39'; Value = '  +0.50%  ' }
    @{ Cell = 'D40'; Value = '0.9029' }
    @{ Cell = 'E40'; Value = '  -0.28%  ' }
    @{ Cell = 'D41'; Value = '1.976' }
    @{ Cell = 'E41'; Value = '  -4.65%  ' }
    @{ Cell = 'D42'; Value = '105.34' }
    @{ Cell = 'E42'; Value = '  -1.39%  ' }
    @{ Cell = 'D43'; Value = '0.9944' }
    @{ Cell = 'D44'; Value = '0.4099' }
    @{ Cell = 'E44'; Value = '  -4.28%  ' }
    @{ Cell = 'D45'; Value = '5.527' }
    @{ Cell = 'E45'; Value = '  -6.42%  ' }
    @{ Cell = 'D46'; Value = '7.036' }
    @{ Cell = 'E46'; Value = '  -5.27%  ' }
    @{ Cell = 'D47'; Value = '61.17' }
    @{ Cell = 'E47'; Value = '  -5.99%  ' }
    @{ Cell = 'D48'; Value = '0.1204' }
    @{ Cell = 'E48'; Value = '  -5.88%  ' }
    @{ Cell = 'D49'; Value = '8.813' }
    @{ Cell = 'E49'; Value = '  -1.41%  ' }
    @{ Cell = 'E50'; Value = '  -5.72%  ' }
    @{ Cell = 'B51'; Value = 'Cronos' }
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro' }
    @{ Cell = 'D51'; Value = '0.05592' }
    @{ Cell = 'E51'; Value = '  -0.77%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $col = $u.Cell -replace '[0-9]+$', ''
    if ($col -eq 'D' -or $col -eq 'E') {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $u.Value
}
